$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear stale forecast values in rows 2-6 (columns C and E)
$ws.Range("E2").ClearContents()

$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()

$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()

$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()

$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Update recomputed forecast values in rows 7-19 (columns C and E)
$ws.Range("C7").Value = 1.785377844167058
$ws.Range("E7").Value = 2.333075171696652

$ws.Range("C8").Value = 5.477304442308206
$ws.Range("E8").Value = 4.052456259163839

$ws.Range("C9").Value = 4.666532690711245
$ws.Range("E9").Value = 3.659383764712709

$ws.Range("C10").Value = 5.266214435142658
$ws.Range("E10").Value = 4.181342739750682

$ws.Range("C11").Value = 4.811826107786477
$ws.Range("E11").Value = 4.131858242365549

$ws.Range("C12").Value = 5.91185619417105
$ws.Range("E12").Value = 4.365509285986957

$ws.Range("C13").Value = 5.114185474093769
$ws.Range("E13").Value = 5.472991335528654

$ws.Range("C14").Value = 2.167530781895133
$ws.Range("E14").Value = 2.573593955528963

$ws.Range("C15").Value = 0.5766229317536675
$ws.Range("E15").Value = 4.059584075094214

$ws.Range("C16").Value = 2.288114387968587
$ws.Range("E16").Value = 3.463553906111505

$ws.Range("C17").Value = -2.013802094285932
$ws.Range("E17").Value = 2.374210810973465

$ws.Range("C18").Value = -0.5865622195987186
$ws.Range("E18").Value = 2.431929210693595

$ws.Range("C19").Value = 0.7174582534189566
$ws.Range("E19").Value = 2.061048937680932
